# The workbook contains a weekly price table for "Tuna" (prickly pear) at
# "Terminal Hortofrutícola Agro Chillán". This edit adds the data for a new
# (most recent) week on top of the table: a new row is inserted at row 8,
# pushing the existing rows 8-17 down to rows 9-18, exactly like the data
# source would when a new week's record is prepended to the series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8; this shifts current rows 8-17 down to 9-18,
# carrying their values/styles with them (matching the diff exactly).
$ws.Rows.Item(8).Insert()

# The columns A, B, C, E, F, G, H, I, J, K, L are constant for every data row
# in this sheet (same market/region/product/category/variety/quality), so we
# copy them from the row directly below (which now holds the former row 8
# data) into the freshly inserted row 8.
$ws.Cells.Item(8, 1).Value = $ws.Cells.Item(9, 1).Value2    # Mercado ID
$ws.Cells.Item(8, 2).Value = $ws.Cells.Item(9, 2).Value2    # Mercado
$ws.Cells.Item(8, 3).Value = $ws.Cells.Item(9, 3).Value2    # Región
$ws.Cells.Item(8, 5).Value = $ws.Cells.Item(9, 5).Value2    # Codreg
$ws.Cells.Item(8, 6).Value = $ws.Cells.Item(9, 6).Value2    # Tipo
$ws.Cells.Item(8, 7).Value = $ws.Cells.Item(9, 7).Value2    # Producto ID
$ws.Cells.Item(8, 8).Value = $ws.Cells.Item(9, 8).Value2    # Producto
$ws.Cells.Item(8, 9).Value = $ws.Cells.Item(9, 9).Value2    # Categoría ID
$ws.Cells.Item(8, 10).Value = $ws.Cells.Item(9, 10).Value2  # Categoría
$ws.Cells.Item(8, 11).Value = $ws.Cells.Item(9, 11).Value2  # Variedad
$ws.Cells.Item(8, 12).Value = $ws.Cells.Item(9, 12).Value2  # Calidad

# New, week-specific values for the inserted row 8.
$ws.Cells.Item(8, 4).Value = 45062                    # Fecha
$ws.Cells.Item(8, 13).Value = 90                      # Volumen
$ws.Cells.Item(8, 14).Value = 13000                   # Precio mínimo
$ws.Cells.Item(8, 15).Value = 14000                   # Precio máximo
$ws.Cells.Item(8, 16).Value = 13444                   # Precio promedio ponderado
$ws.Cells.Item(8, 17).Value = "$/caja 18 kilos"       # Unidad de comercialización
$ws.Cells.Item(8, 18).Value = "Región Metropolitana"  # Origen
$ws.Cells.Item(8, 19).Value = 747                     # Precio $/Kg
$ws.Cells.Item(8, 20).Value = 18                      # Kg / unidad
